# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''19.929.67'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -8.13%  '

# Row 3
$ws.Range("D3").Value = '''1.406.01'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -8.29%  '

# Row 4
$ws.Range("D4").Value = '''1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.16%  '

# Row 5
$ws.Range("D5").Value = '''1.002'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.09%  '

# Row 6
$ws.Range("D6").Value = '''271.18'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.10%  '

# Row 7
$ws.Range("D7").Value = '''0.3685'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -5.95%  '

# Row 8
$ws.Range("D8").Value = '''0.3054'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.32%  '

# Row 9
$ws.Range("D9").Value = '''38.92'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -8.00%  '

# Row 10
$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").Value = '''0.9827'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.00%  '

# Row 11
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = '''0.06490'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -9.47%  '

# Row 12
$ws.Range("D12").Value = '''1.003'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.19%  '

# Row 13
$ws.Range("D13").Value = '''5.295'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.96%  '

# Row 14
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = '''6.093'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -7.55%  '

# Row 15
$ws.Range("B15").Value = 'Solana'
$ws.Range("C15").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D15").Value = '''16.72'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -9.74%  '

# Row 16
$ws.Range("D16").Value = '''1.408.96'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -8.41%  '

# Row 17
$ws.Range("D17").Value = '''0.000009998'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -8.85%  '

# Row 18
$ws.Range("D18").Value = '''0.05706'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -13.34%  '

# Row 19
$ws.Range("D19").Value = '''72.49'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -12.67%  '

# Row 20
$ws.Range("E20").Value = '  +0.07%  '

# Row 21
$ws.Range("D21").Value = '''5.530'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -9.51%  '

# Row 22
$ws.Range("D22").Value = '''14.24'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -7.57%  '

# Row 23
$ws.Range("D23").Value = '''10.75'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.76%  '

# Row 24
$ws.Range("D24").Value = '''2.273'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.25%  '

# Row 25
$ws.Range("D25").Value = '''19.953.28'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -8.05%  '

# Row 26
$ws.Range("D26").Value = '''2.193'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.74%  '

# Row 27
$ws.Range("D27").Value = '''137.42'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.81%  '

# Row 28
$ws.Range("D28").Value = '''16.58'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -9.52%  '

# Row 29
$ws.Range("D29").Value = '''1.566.13'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -8.60%  '

# Row 30
$ws.Range("D30").Value = '''107.95'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.89%  '

# Row 31
$ws.Range("D31").Value = '''3.849'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -20.52%  '

# Row 32
$ws.Range("D32").Value = '''5.205'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -11.50%  '

# Row 33
$ws.Range("D33").Value = '''0.8048'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -16.02%  '

# Row 34
$ws.Range("D34").Value = '''0.07657'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.95%  '

# Row 35
$ws.Range("D35").Value = '''8.362'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.19%  '

# Row 36
$ws.Range("D36").Value = '''0.05805'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.29%  '

# Row 37
$ws.Range("D37").Value = '''1.001'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.01%  '

# Row 38
$ws.Range("D38").Value = '''4.734'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.46%  '

# Row 39
$ws.Range("E39").Value = '  -4.10%  '

# Row 40
$ws.Range("D40").Value = '''0.02019'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -8.31%  '

# Row 41
$ws.Range("D41").Value = '''10.07'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -7.51%  '

# Row 42
$ws.Range("B42").Value = 'WEMIXTOKEN'
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D42").Value = '''1.290'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -11.16%  '

# Row 43
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '''1.056'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -10.42%  '

# Row 44
$ws.Range("D44").Value = '''0.5239'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -8.98%  '

# Row 45
$ws.Range("D45").Value = '''3.508'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.95%  '

# Row 46
$ws.Range("D46").Value = '''12.05'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.30%  '

# Row 47
$ws.Range("D47").Value = '''0.5063'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.75%  '

# Row 48
$ws.Range("D48").Value = '''1.785'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.32%  '

# Row 49
$ws.Range("D49").Value = '''109.36'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.50%  '

# Row 50
$ws.Range("D50").Value = '''1.033'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -10.35%  '

# Row 51
$ws.Range("D51").Value = '''1.001'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.02%  '
